# Update gh-pages "南宁-漫展信息" workbook data: bump several "want-to-go"
# counters and insert a new exhibition row ("南宁·0713国乙ONLY") into both
# the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    # Force plain text so date-looking strings ("2024-07-13") are not
    # auto-converted into date serials by Excel's input parser.
    $ws.Cells.Item($row, $col).NumberFormat = "General"
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

function Set-NumberCell($ws, $row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

function Insert-NewEventRow($ws, $newRow) {
    # Insert a blank row, shifting the rest of the table down, then copy the
    # (bold + bordered) formatting of column A from the row below onto the
    # freshly inserted index cell.
    $ws.Rows.Item($newRow).Insert()
    $ws.Cells.Item($newRow + 1, 1).Copy()
    $ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
}

function Fill-EventRow($ws, $row, $date, $name, $place, $timeRange, $want, $price, $link, $cover) {
    Set-TextCell $ws $row 2 $date
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 4).Value = $place
    Set-TextCell $ws $row 5 $timeRange
    Set-NumberCell $ws $row 6 $want
    Set-NumberCell $ws $row 7 $price
    $ws.Cells.Item($row, 8).Value = $link
    $ws.Cells.Item($row, 9).Value = $cover
}

function Renumber-IndexColumn($ws, $lastRow) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).NumberFormat = "General"
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

# ---------------------------------------------------------------------
# Sheet "展览": new row goes in at row 10 (between "首届童话梦境Lolita茶会"
# and "广西·首届明日方舟only展"), table grows from 14 to 15 data rows.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

Set-NumberCell $ws1 3 6 3922   # 南宁·AP动漫游戏嘉年华
Set-NumberCell $ws1 4 6 2314   # 南宁·布谷鸟动漫展4th
Set-NumberCell $ws1 6 6 14     # 宾阳·荷止国风动漫展
Set-NumberCell $ws1 8 6 189    # 南宁·小蜜蜂动漫嘉年华2.0
Set-NumberCell $ws1 9 6 111    # 南宁·首届童话梦境Lolita茶会

Insert-NewEventRow $ws1 10
Fill-EventRow $ws1 10 "2024-07-13" "南宁·0713国乙ONLY" "亭洪路45号 水明漾宴会中心" `
    "2024.07.13 09:30-07.13 21:00" 2 68 `
    "https://show.bilibili.com/platform/detail.html?id=86378" `
    "//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg"

# Rows after the insert point keep their old values (shifted, but not
# renumbered) and the two "want-to-go" counters below also moved/changed.
Set-NumberCell $ws1 12 6 1454  # 南宁·AB动漫游戏嘉年华 (was row 11)
Set-NumberCell $ws1 13 6 258   # 横州·第二届海棠动漫游戏嘉年华 (was row 12)
Set-NumberCell $ws1 14 6 2616  # 良牙夏典 (was row 13)
Set-NumberCell $ws1 15 6 181   # 南宁·蔚蓝档案only (was row 14)

Renumber-IndexColumn $ws1 15

# ---------------------------------------------------------------------
# Sheet "全部类型": new row goes in at row 11 (between "首届童话梦境Lolita
# 茶会" and "广西·首届明日方舟only展"), table grows from 17 to 18 data rows.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

Set-NumberCell $ws4 3 6 3922   # 南宁·AP动漫游戏嘉年华
Set-NumberCell $ws4 4 6 2314   # 南宁·布谷鸟动漫展4th
Set-NumberCell $ws4 6 6 14     # 宾阳·荷止国风动漫展
Set-NumberCell $ws4 9 6 189    # 南宁·小蜜蜂动漫嘉年华2.0
Set-NumberCell $ws4 10 6 111   # 南宁·首届童话梦境Lolita茶会

Insert-NewEventRow $ws4 11
Fill-EventRow $ws4 11 "2024-07-13" "南宁·0713国乙ONLY" "亭洪路45号 水明漾宴会中心" `
    "2024.07.13 09:30-07.13 21:00" 2 68 `
    "https://show.bilibili.com/platform/detail.html?id=86378" `
    "//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg"

Set-NumberCell $ws4 15 6 1454  # 南宁·AB动漫游戏嘉年华 (was row 14)
Set-NumberCell $ws4 16 6 258   # 横州·第二届海棠动漫游戏嘉年华 (was row 15)
Set-NumberCell $ws4 17 6 2616  # 良牙夏典 (was row 16)
Set-NumberCell $ws4 18 6 181   # 南宁·蔚蓝档案only (was row 17)

Renumber-IndexColumn $ws4 18

Write-Host "done"
